# Working With Data.xlsx - reformat the single JSON-ish "questions" blob
# that lives in A2 (the diff's sharedStrings.xml shows the same content,
# just pretty-printed instead of single-line) and relocate it to A1,
# dropping the old numeric placeholder (A1 = 0) and its bordered/bold
# "header" style so the sheet ends up as a single, unformatted A1 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are presented with the table shown below, in which all the information is accurate.Which column requires attention before you can analyze the data?NameClassDate of BirthTest ResultJames4PW30th November80.5Mohammed4A5/2/1283Sofia4A1st Dec 201167George4PW3/8/1276.5Nicola4PW5/4/1278Lena4G24-Jan-1292",
        "ques_type": 2,
        "options": [
            "Name",
            "Class",
            "Date of Birth",
            "Test Result"
        ],
        "score": "Date of Birth"
    },
    {
        "title": "How would you find the mean of a set of numbers?",
        "ques_type": 2,
        "options": [
            "Find the number that appears most frequently.",
            "Calculate the sum of the numbers and divide by the amount of numbers in the set.",
            "Find the number halfway between the smallest number and the largest number.",
            "Put all the numbers in ascending order and find the number in the middle."
        ],
        "score": "Calculate the sum of the numbers and divide by the amount of numbers in the set."
    },
    {
        "title": "You are working as a marketing analyst for an ice cream company, and you are presented with data from a survey on people\u2019s favorite ice cream flavors. In the survey, people were asked to select their favorite flavor from a list of 25 options, and over 800 people responded. Your manager has asked you to produce a quick chart to illustrate and compare the popularity of all the flavors.  Which type of chart would be best suited to the task?",
        "ques_type": 2,
        "options": [
            "Scatter plot",
            "Pie chart",
            "Bar chart",
            "Line chart"
        ],
        "score": "Pie chart"
    },
    {
        "title": "You work for a furniture sales company that has stores across the city. You are presented with the chart shown below comparing sales in two different regions of the city.Which of the following statements can be concluded from the information in the chart?",
        "ques_type": 15,
        "options": [
            "Sales generally decreased over the six months shown.",
            "More items were sold in the North Region than the South Region during the six months shown.",
            "More items will be sold in July than were sold in June.",
            "January was the month with the fewest sales.",
            "The average number of items sold per month is lower in the South region",
            "The North region makes a higher profit than the South region"
        ],
        "score": [
            "More items were sold in the North Region than the South Region during the six months shown.",
            "January was the month with the fewest sales.",
            "The average number of items sold per month is lower in the South region"
        ]
    }
]
'@

# A1 currently holds the placeholder number 0 styled with a thin border +
# bold, centered font; A2 holds the real text. Clear both out...
$ws.Range("A1").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A2").ClearContents()

# ...then write the reformatted text into A1 using the sheet's default
# (unstyled) formatting, matching the new single-row layout.
$ws.Range("A1").Value = $text

# Setting a multi-line value auto-expands the row height; AutoFit restores
# it to the sheet's default so row 1 keeps no explicit height override.
$ws.Rows.Item(1).AutoFit()
